$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the A-column group labels first, in the order the new distinct
# values are first encountered (VAD, Heart transplant, Organ donor), so
# new shared-string entries get appended in that order - matching the
# target workbook's sharedStrings table layout.
$ws.Range("A5").Value = "VAD"
$ws.Range("A6").Value = "VAD"
$ws.Range("A7").Value = "VAD"

$ws.Range("A8").Value = "Heart transplant"
$ws.Range("A9").Value = "Heart transplant"
$ws.Range("A10").Value = "Heart transplant"
$ws.Range("A11").Value = "Heart transplant"
$ws.Range("A12").Value = "Heart transplant"

$ws.Range("A2").Value = "Organ donor"
$ws.Range("A3").Value = "Organ donor"
$ws.Range("A4").Value = "Organ donor"

# Header row last (new strings factor_1, grouping appended after the
# group-label strings above; "y" already exists so it's reused as-is).
$ws.Range("A1").Value = "factor_1"
$ws.Range("B1").Value = "grouping"
$ws.Range("C1").Value = "y"

# B-column tag/subgroup labels reuse existing shared strings (a..g).
$ws.Range("B2").Value = "a"
$ws.Range("B3").Value = "a"
$ws.Range("B4").Value = "b"
$ws.Range("B5").Value = "c"
$ws.Range("B6").Value = "d"
$ws.Range("B7").Value = "e"
$ws.Range("B8").Value = "c"
$ws.Range("B9").Value = "c"
$ws.Range("B10").Value = "d"
$ws.Range("B11").Value = "f"
$ws.Range("B12").Value = "g"

# Column A width to fit the longest label ("Heart transplant")
$ws.Columns.Item(1).ColumnWidth = 14.6

# Selection moves to B2
$ws.Range("B2").Select()
